$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended after the 2025-09-06 run.
# Force the date column to be stored as literal text (matching the
# existing rows above, which are plain strings rather than real dates),
# then drop the temporary "@" number format so the cell's style matches
# its un-styled neighbours (A2:A4).
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "09/06/2025"
$ws.Range("A5").ClearFormats()

$ws.Range("B5").Value = 0.1259940523634941
$ws.Range("C5").Value = 0.8740059476365059
